$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# STAGE sheet: latest automation run data gets written into row 2-5
# ------------------------------------------------------------------
$stage = $wb.Worksheets.Item("STAGE")

$stage.Range("A2").Value = "ECLocation65237"
$stage.Range("D2").Value = "admin"
$stage.Range("E2").Value = "controller"
$stage.Range("H2").Value = "FPK12Exam54223"
$stage.Range("I2").Value = "FPK12Schedule48334"

$stage.Range("E3").Value = "48196"
$stage.Range("E4").Value = "51072"
$stage.Range("E5").Value = "68891"

# ------------------------------------------------------------------
# View state: LMSProd's own lingering selection moves to E8, then
# STAGE becomes the active/selected tab with selection on E2.
# ------------------------------------------------------------------
$lms = $wb.Worksheets.Item("LMSProd")
$lms.Activate()
$lms.Range("E8").Select()

$stage.Activate()
$stage.Range("E2").Select()
